$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-14T15:34:52+00:00"
$meta.Range("B12").Value = "Dose d'antigène"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("M2").Value = "Dose d'antigène"
$elements.Range("L4").Value = "Dose d'antigène"
$elements.Range("M4").Value = "Dose d'antigène"
